$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1, matching the formatting of the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H7 with 0 (plain numeric values, no special style)
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
